# RELEASE 1.0, drive cache removed
# Remove the "drive cache" column (C) from the Projects sheet, and rename
# the "Persons" sheet headers/column to reflect the new tracking scheme.

$wb = $excel.ActiveWorkbook

# --- Sheet "Projects" (sheet1): drop the "Примечание" column (C) ---
$ws1 = $wb.Worksheets.Item(1)

# Clear the header cell (keeps its bold header style) and the three data
# cells below it.
$ws1.Range("C1").ClearContents()
$ws1.Range("C2:C4").ClearContents()

# --- Sheet "Persons" (sheet2): rename headers for the new workflow ---
$ws2 = $wb.Worksheets.Item(2)

$ws2.Range("C1").Value = "ТрЗ за смену статуса"
$ws2.Range("A1").Value = "ФИО"

# Widen column C to (best) fit the new, longer header text.
$ws2.Columns.Item(3).ColumnWidth = 17.9

# Move the saved selection to C10, as in the authored workbook.
$null = $ws2.Range("C10").Select()

# Re-select sheet1's saved cell (B12) last, so Projects stays the active
# tab, matching the authored workbook.
$null = $ws1.Range("B12").Select()
